$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

$ws.Range("D47").Select()
